# Update countries & provincias Spain
# Applies the 19-Aug-2020 data refresh to the "Pais" sheet:
#   - Updates the "last updated" timestamp (A1)
#   - Refreshes case counts for several countries (rows keep their rank,
#     only the B..H figures move)
#   - Two country pairs swap rank (their row position in the table),
#     which also swaps which country label + data sits on each row:
#       row 173/174 : Islas Feroe <-> Birmania
#       row 213/214 : Montserrat  <-> Islas Malvinas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Timestamp header -----------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 19 de Agosto de 2020 a las 12:10"

# ---- Straight numeric refreshes (country stays on its row) -----------
# Estados Unidos (row 4)
$ws.Range("E4").Value = 2469535
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 175092

# Oman (row 38)
$ws.Range("B38").Value = 83606
$ws.Range("C38").Value = 188
$ws.Range("D38").Value = 78188
$ws.Range("E38").Value = 4815
$ws.Range("G38").Value = 6
$ws.Range("H38").Value = 603

# Rumania (row 42)
$ws.Range("B42").Value = 73617
$ws.Range("C42").Value = 1409
$ws.Range("D42").Value = 33566
$ws.Range("E42").Value = 36945
$ws.Range("G42").Value = 32
$ws.Range("H42").Value = 3106

# El Salvador (row 73)
$ws.Range("B73").Value = 23717
$ws.Range("C73").Value = 255
$ws.Range("D73").Value = 11388
$ws.Range("E73").Value = 11696
$ws.Range("G73").Value = 8
$ws.Range("H73").Value = 633

# Consejo Danes para los Refugiados (row 90)
$ws.Range("B90").Value = 9741
$ws.Range("C90").Value = 20
$ws.Range("D90").Value = 8895
$ws.Range("E90").Value = 600
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 246

# Finlandia (row 98)
$ws.Range("B98").Value = 7805
$ws.Range("C98").Value = 29
$ws.Range("E98").Value = 421

# Santo Tome y Principe (row 162)
$ws.Range("D162").Value = 826
$ws.Range("E162").Value = 44

# Timor Oriental (row 203)
$ws.Range("D203").Value = 25
$ws.Range("E203").Value = 0

# ---- Rank swaps: country label + full data row change places ---------
# Islas Feroe (was row 173) <-> Birmania (was row 174)
$ws.Range("A173").Value = "Birmania"
$ws.Range("B173").Value = 379
$ws.Range("C173").Value = 3
$ws.Range("D173").Value = 331
$ws.Range("E173").Value = 42
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 6

$ws.Range("A174").Value = "Islas Feroe"
$ws.Range("B174").Value = 377
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 235
$ws.Range("E174").Value = 142
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

# Montserrat (was row 213) <-> Islas Malvinas (was row 214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 13
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
